{"js": "// Apply the benchmark-table value updates described by the diff.\n// The document body contains a single table with one column; each row\n// holds one statistic. A handful of rows get their text value swapped,\n// and three rows that used to hold a full tab-separated detail line get\n// collapsed down to just their first (summary) value.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"5062\",\n  5: \"0.05147\",\n  6: \"0.00781\",\n  7: \"0.00378\",\n  8: \"0.02550\",\n  9: \"0.02701\",\n  10: \"0.04549\",\n  11: \"4.29034\",\n  43: \"98.4\",\n  44: \"4.29\",\n  45: \"267\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(parseInt(rowIndex, 10), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table value updates described by the diff.\n# The document body contains a single table with one column; each row\n# holds one statistic. A handful of rows get their text value swapped,\n# and three rows that used to hold a full tab-separated detail line get\n# collapsed down to just their first (summary) value.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row number -> new cell text (Word COM tables/cells are 1-indexed).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"5062\"\n    6  = \"0.05147\"\n    7  = \"0.00781\"\n    8  = \"0.00378\"\n    9  = \"0.02550\"\n    10 = \"0.02701\"\n    11 = \"0.04549\"\n    12 = \"4.29034\"\n    44 = \"98.4\"\n    45 = \"4.29\"\n    46 = \"267\"\n}\n\nforeach ($row in $updates.Keys) {\n    $t.Cell($row, 1).Range.Text = $updates[$row]\n}\n"}
